$d = $word.ActiveDocument

# wdReplaceOne = 1 (replace only the current match, not all)
$wdReplaceOne = 1
# wdFindContinue = 1 (do not wrap around when searching)
$wdFindContinue = 1
# wdCollapseEnd = 0
$wdCollapseEnd = 0

function Replace-NextOccurrence($range, $findText, $replaceText) {
    $range.Find.Execute($findText, $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceOne) | Out-Null
    $range.Collapse($wdCollapseEnd)
}

# --- Main body: "A TERE," -> "A QWER," (bold run) ---
$bodyRange = $d.Content
$bodyRange.Find.Execute("TERE", $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, "QWER", $wdReplaceOne) | Out-Null

# --- Header: sequential, ordered replacements matching document order ---
$hdr = $d.Sections(1).Headers(1)
$hdrRange = $hdr.Range
$hdrRange.End = $hdr.Range.StoryLength - 1

Replace-NextOccurrence $hdrRange "TRE"  "QWER"
Replace-NextOccurrence $hdrRange "TERE" "QWER"
Replace-NextOccurrence $hdrRange "Tre"  "Qwer"
Replace-NextOccurrence $hdrRange "Tre"  "Qwer"
Replace-NextOccurrence $hdrRange "Tre"  "Qewr"
Replace-NextOccurrence $hdrRange "Tre"  "Qewr"
Replace-NextOccurrence $hdrRange "Tre"  "Qwer"
Replace-NextOccurrence $hdrRange "tre"  "qwer"
Replace-NextOccurrence $hdrRange "tre"  "qwer"
Replace-NextOccurrence $hdrRange "tre"  "qwer"
